# Apply the text edits described by the commit diff:
#  - Slide 2, "Content Placeholder 2": "diagnosis" -> "assessment"
#  - Slide 3, "Title 1": "Diagnosis" -> "assessment"
#
# Setting TextRange/Paragraph .Text directly to the final string causes the
# interop layer to diff old vs. new text and split the run around the
# changed word (it keeps the unchanged prefix/suffix as separate runs).
# The canonical edit only ever touches the single existing <a:r>, so we
# first stomp the paragraph text with an unrelated placeholder (no shared
# prefix/suffix with either the old or the new text) and then set the real
# text; since neither step shares overlap with its predecessor, the whole
# paragraph is rewritten as one run, matching the target XML exactly.

$p = $ppt.ActivePresentation

$s2 = $p.Slides.Item(2)
$shGoals = $s2.Shapes.Item(2)
$para1 = $shGoals.TextFrame.TextRange.Paragraphs(1)
$para1.Text = "zzz_placeholder_zzz"
$para1.Text = "Quick + Accurate assessment and treatment of the patient"

$s3 = $p.Slides.Item(3)
$shTitle = $s3.Shapes.Item(1)
$tr3 = $shTitle.TextFrame.TextRange
$tr3.Text = "zzz_placeholder_zzz"
$tr3.Text = "Quick and Accurate assessment and Treatment"
